# Sat, May 30, 2020  3:04:42 PM
#
# The deck currently uses the custom "Integral" theme (dark green /
# lime accents) for the slide master / slides, while the plain
# built-in "Office Theme" palette sits unused on the notes master.
# This edit swaps the two: the slide master's theme is recoloured to
# the standard "Office Theme" palette (what ships in ppt/theme/theme2.xml
# becomes the Office colours), matching the effect of the author
# re-applying the default Office colour scheme to the design.
#
# PowerPoint's ThemeColorScheme exposes the 12 theme colour slots as
# an OLE_COLOR (0x00BBGGRR) on ThemeColor.RGB -- there's no RGB()
# helper in this shell, so build the packed integer by hand.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Target palette = standard Office Theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), in ThemeColorScheme.Colors() index order.
$colors.Colors(1).RGB  = RGBVal 0   0   0    # dk1      000000
$colors.Colors(2).RGB  = RGBVal 255 255 255  # lt1      FFFFFF
$colors.Colors(3).RGB  = RGBVal 68  84  106  # dk2      44546A
$colors.Colors(4).RGB  = RGBVal 231 230 230  # lt2      E7E6E6
$colors.Colors(5).RGB  = RGBVal 91  155 213  # accent1  5B9BD5
$colors.Colors(6).RGB  = RGBVal 237 125 49   # accent2  ED7D31
$colors.Colors(7).RGB  = RGBVal 165 165 165  # accent3  A5A5A5
$colors.Colors(8).RGB  = RGBVal 255 192 0    # accent4  FFC000
$colors.Colors(9).RGB  = RGBVal 68  114 196  # accent5  4472C4
$colors.Colors(10).RGB = RGBVal 112 173 71   # accent6  70AD47
$colors.Colors(11).RGB = RGBVal 5   99  193  # hlink    0563C1
$colors.Colors(12).RGB = RGBVal 149 79  114  # folHlink 954F72

# Best-effort: also try to relabel the design / colour scheme from
# "Integral" to "Office Theme" so the palette's name matches its new
# values (no-op on hosts that treat theme naming as read-only).
$design.Name = "Office Theme"
$theme.Name = "Office Theme"
$colors.Name = "Office"
